$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 30,3
$data[0,0] = -9.514568328857422
$data[0,1] = -7.287443161010742
$data[0,2] = -4.99766206741333
$data[1,0] = -7.823380947113037
$data[1,1] = -6.074180126190186
$data[1,2] = -4.752357959747314
$data[2,0] = -3.637385606765747
$data[2,1] = -2.725360631942749
$data[2,2] = -3.621543645858765
$data[3,0] = -10.94634628295898
$data[3,1] = -2.156153917312622
$data[3,2] = -9.503169059753418
$data[4,0] = 66.69537353515625
$data[4,1] = -35.59264755249023
$data[4,2] = -8.555927276611328
$data[5,0] = -6.669784069061279
$data[5,1] = -7.415677547454834
$data[5,2] = 10.31031608581543
$data[6,0] = -6.682662963867188
$data[6,1] = 0.0761735439300537
$data[6,2] = -2.404594898223877
$data[7,0] = -4.13407564163208
$data[7,1] = -28.60597419738769
$data[7,2] = 8.249073028564453
$data[8,0] = -12.32790374755859
$data[8,1] = 6.361005783081055
$data[8,2] = -18.50937080383301
$data[9,0] = -0.0332281589508056
$data[9,1] = -8.605781555175781
$data[9,2] = -5.739476203918457
$data[10,0] = -14.07493591308594
$data[10,1] = -31.83533477783203
$data[10,2] = -4.278344631195068
$data[11,0] = 2.169375419616699
$data[11,1] = 6.375825881958008
$data[11,2] = 16.93547058105469
$data[12,0] = -3.887731313705444
$data[12,1] = 1.673339605331421
$data[12,2] = 2.357208251953125
$data[13,0] = 3.855255126953125
$data[13,1] = -33.85980224609375
$data[13,2] = 3.360419273376465
$data[14,0] = 30.62849044799805
$data[14,1] = 7.678761005401611
$data[14,2] = -9.237998962402344
$data[15,0] = -11.34725379943848
$data[15,1] = -16.60527801513672
$data[15,2] = -12.14533615112305
$data[16,0] = 6.219323635101318
$data[16,1] = -10.72451782226562
$data[16,2] = 26.53547286987305
$data[17,0] = -4.500537395477295
$data[17,1] = 9.60122776031494
$data[17,2] = -3.719542026519776
$data[18,0] = -59.25642776489258
$data[18,1] = -72.75296783447266
$data[18,2] = 58.0263671875
$data[19,0] = 43.50658416748047
$data[19,1] = 8.478635787963867
$data[19,2] = -37.43244552612305
$data[20,0] = -19.30278778076172
$data[20,1] = -6.771676063537598
$data[20,2] = -17.75639343261719
$data[21,0] = -18.29881477355957
$data[21,1] = -37.5744743347168
$data[21,2] = 5.842066764831543
$data[22,0] = -2.273155212402344
$data[22,1] = 8.59691047668457
$data[22,2] = -6.313179969787598
$data[23,0] = -0.403256893157959
$data[23,1] = 4.687671661376953
$data[23,2] = -1.856612205505371
$data[24,0] = 4.316394805908203
$data[24,1] = -26.35572052001953
$data[24,2] = -17.98580360412598
$data[25,0] = 31.11298370361328
$data[25,1] = 3.278896331787109
$data[25,2] = -11.80455207824707
$data[26,0] = -10.42159271240234
$data[26,1] = -19.3218994140625
$data[26,2] = -14.07432746887207
$data[27,0] = 6.261336803436279
$data[27,1] = -9.548392295837402
$data[27,2] = 26.6091365814209
$data[28,0] = -5.242365837097168
$data[28,1] = 18.12157821655273
$data[28,2] = 3.32082748413086
$data[29,0] = -17.76671600341797
$data[29,1] = -38.77998352050781
$data[29,2] = 33.62932968139648

$ws.Range("A2:C31").Value = $data
